$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    5  = -2
    6  = 4
    7  = 0
    8  = -1
    9  = -7
    10 = -1
    11 = 2
    13 = -3
    14 = 1
    15 = 2
    16 = 11
    17 = 3
    18 = -2
    20 = 5
    22 = 2
    23 = 8
    25 = -5
    26 = -3
    27 = 3
    28 = -1
    29 = 1
    30 = -2
    31 = 3
    32 = -7
    33 = -2
    34 = -2
    36 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
